# Weekly update: new "Clementina" price report (week of 2022-07-05) is inserted
# as rows 113-114, pushing the previous "Murcott" rows (week of 2021-09-03) down
# to rows 115-116.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 113:114, shifting the old rows 113-114 down to 115-116.
$ws.Rows("113:114").Insert()

# Row 113 - Clementina, Primera
$ws.Range("A113").Value = 11
$ws.Range("B113").Value = "Vega Monumental Concepción"
$ws.Range("C113").Value = "Bíobío"
$ws.Range("D113").Value = 44747
$ws.Range("E113").Value = 8
$ws.Range("F113").Value = "Fruta"
$ws.Range("G113").Value = 100102
$ws.Range("H113").Value = "Cítricos"
$ws.Range("I113").Value = 100102004
$ws.Range("J113").Value = "Mandarina"
$ws.Range("K113").Value = "Clementina"
$ws.Range("L113").Value = "Primera"
$ws.Range("M113").Value = 100
$ws.Range("N113").Value = 7000
$ws.Range("O113").Value = 8000
$ws.Range("P113").Value = 7500
$ws.Range("Q113").Value = "`$/bandeja 18 kilos"
$ws.Range("R113").Value = "Región de O'Higgins"
$ws.Range("S113").Value = 417
$ws.Range("T113").Value = 18

# Row 114 - Clementina, Segunda
$ws.Range("A114").Value = 11
$ws.Range("B114").Value = "Vega Monumental Concepción"
$ws.Range("C114").Value = "Bíobío"
$ws.Range("D114").Value = 44747
$ws.Range("E114").Value = 8
$ws.Range("F114").Value = "Fruta"
$ws.Range("G114").Value = 100102
$ws.Range("H114").Value = "Cítricos"
$ws.Range("I114").Value = 100102004
$ws.Range("J114").Value = "Mandarina"
$ws.Range("K114").Value = "Clementina"
$ws.Range("L114").Value = "Segunda"
$ws.Range("M114").Value = 50
$ws.Range("N114").Value = 6000
$ws.Range("O114").Value = 6000
$ws.Range("P114").Value = 6000
$ws.Range("Q114").Value = "`$/bandeja 18 kilos"
$ws.Range("R114").Value = "Región de O'Higgins"
$ws.Range("S114").Value = 333
$ws.Range("T114").Value = 18
